$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the "Upute za pokretanje" section so we only touch that part
# of the document (similar phrases appear earlier too, describing the
# already-working game).
# ------------------------------------------------------------------
$anchor = $d.Content
$foundAnchor = $anchor.Find.Execute("Upute za pokretanje", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$sectionStart = $anchor.End

# 1) Remove the run " nakon pozicioniranja u direktorij u kojem se nalazi snake.py file"
$r1 = $d.Range($sectionStart, $d.Content.End)
$found1 = $r1.Find.Execute(" nakon pozicioniranja u direktorij u kojem se nalazi snake.py file", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $r1.Delete()
}

# 2) Replace "slijedećim naredbama:" with the new trailing explanation text
$r2 = $d.Range($sectionStart, $d.Content.End)
$found2 = $r2.Find.Execute("slijedećim naredbama:", $false, $false, $false, $false, $false, $true, 1, $false, "naredbom uz pretpostavku da je instaliranja pygame knjižnice", 2)

# 3) Delete the three paragraphs that used to describe the venv activation
#    steps ("1.ucitavanja virtualnog okruženja", ".venv\Scripts\activate",
#    "2. pokretanje skripte"), merging straight into the "python snake.py"
#    paragraph that follows them.
$r3 = $d.Range($sectionStart, $d.Content.End)
$found3a = $r3.Find.Execute("1.ucitavanja virtualnog", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3a) {
    $delStart = $r3.Start

    $r4 = $d.Range($sectionStart, $d.Content.End)
    $found3b = $r4.Find.Execute("2. pokretanje skripte", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found3b) {
        # Grow by one character so the trailing paragraph mark is included
        # too - otherwise Delete() on a ragged, multi-paragraph range is a
        # silent no-op.
        [void]$r4.MoveEnd(1, 1)
        $delEnd = $r4.End

        $rdel = $d.Range($delStart, $delEnd)
        $rdel.Delete()
    }
}

# 4) Remove the first-line indent that used to set the "python snake.py"
#    line apart from the (now deleted) numbered steps above it.
$r5 = $d.Range($sectionStart, $d.Content.End)
$found4 = $r5.Find.Execute("python snake.py", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found4) {
    $r5.ParagraphFormat.FirstLineIndent = 0
}
